# Update gh-pages to output generated at 456a3b4
# Updates the "想去人数" (want-to-go count) column F on the 展览, 演出
# and 全部类型 sheets to reflect refreshed data.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1580
$ws1.Range("F5").Value = 277
$ws1.Range("F6").Value = 63
$ws1.Range("F7").Value = 1658
$ws1.Range("F8").Value = 10229
$ws1.Range("F9").Value = 173
$ws1.Range("F14").Value = 7086
$ws1.Range("F16").Value = 665
$ws1.Range("F17").Value = 41
$ws1.Range("F18").Value = 56
$ws1.Range("F19").Value = 236

# --- Sheet "演出" (sheet2) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 556
$ws2.Range("F4").Value = 1

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1580
$ws4.Range("F5").Value = 277
$ws4.Range("F7").Value = 63
$ws4.Range("F8").Value = 1658
$ws4.Range("F9").Value = 556
$ws4.Range("F10").Value = 1
$ws4.Range("F11").Value = 10229
$ws4.Range("F12").Value = 173
$ws4.Range("F17").Value = 7086
$ws4.Range("F19").Value = 665
$ws4.Range("F20").Value = 41
$ws4.Range("F21").Value = 56
$ws4.Range("F22").Value = 236
